# Auto-generated Excel COM-interop script to apply profit-sheet value updates
# per the commit diff (scheduled runner market-price refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1386.4
$ws.Range("I19").Value = 1235.7142
$ws.Range("J19").Value = 1518.25
$ws.Range("K19").Value = 1235.7142
$ws.Range("L19").Value = 1518.25
$ws.Range("M19").Value = -1060.7142
$ws.Range("N19").Value = -1868.25
$ws.Range("H40").Value = 2592.8333
$ws.Range("I40").Value = 2403.4285
$ws.Range("J40").Value = 2858
$ws.Range("K40").Value = 2403.4285
$ws.Range("L40").Value = 2858
$ws.Range("M40").Value = -2228.4285
$ws.Range("N40").Value = -3208
$ws.Range("H51").Value = 9859.5
$ws.Range("J51").Value = 9946.166999999999
$ws.Range("L51").Value = 9946.166999999999
$ws.Range("N51").Value = -10914.167
$ws.Range("H87").Value = 30000
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37496
$ws.Range("H90").Value = 30000
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -117480
$ws.Range("H99").Value = 419.17648
$ws.Range("J99").Value = 617.25
$ws.Range("L99").Value = 1851.75
$ws.Range("N99").Value = -4847.75
$ws.Range("H113").Value = 6513.647
$ws.Range("I113").Value = 5495.7144
$ws.Range("J113").Value = 7226.2
$ws.Range("K113").Value = 5495.7144
$ws.Range("L113").Value = 7226.2
$ws.Range("M113").Value = -2241.7144
$ws.Range("N113").Value = -13734.2
$ws.Range("H116").Value = 5908.8184
$ws.Range("I116").Value = 6130
$ws.Range("J116").Value = 5643.4
$ws.Range("K116").Value = 6130
$ws.Range("L116").Value = 5643.4
$ws.Range("M116").Value = -2688
$ws.Range("N116").Value = -12527.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1160.2307
$ws.Range("I2").Value = 1061.25
$ws.Range("J2").Value = 1318.6
$ws.Range("K2").Value = 1061.25
$ws.Range("L2").Value = 1318.6
$ws.Range("M2").Value = -948.25
$ws.Range("N2").Value = -1544.6
$ws.Range("H61").Value = 2380.125
$ws.Range("I61").Value = 2265.6667
$ws.Range("K61").Value = 2265.6667
$ws.Range("M61").Value = -2053.6667
$ws.Range("H116").Value = 1160.2307
$ws.Range("I116").Value = 1061.25
$ws.Range("J116").Value = 1318.6
$ws.Range("K116").Value = 1061.25
$ws.Range("L116").Value = 1318.6
$ws.Range("M116").Value = 1232.75
$ws.Range("N116").Value = -5906.6
$ws.Range("H136").Value = 2380.125
$ws.Range("I136").Value = 2265.6667
$ws.Range("K136").Value = 6797.000100000001
$ws.Range("M136").Value = -4247.000100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1160.2307
$ws.Range("I3").Value = 1061.25
$ws.Range("J3").Value = 1318.6
$ws.Range("K3").Value = 1061.25
$ws.Range("L3").Value = 1318.6
$ws.Range("M3").Value = -947.25
$ws.Range("N3").Value = -1546.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 482.33334
$ws.Range("I16").Value = 479
$ws.Range("J16").Value = 489
$ws.Range("K16").Value = 479
$ws.Range("L16").Value = 489
$ws.Range("M16").Value = -192
$ws.Range("N16").Value = -1063
$ws.Range("H107").Value = 1630.2307
$ws.Range("I107").Value = 1609.4
$ws.Range("J107").Value = 1699.6666
$ws.Range("K107").Value = 1609.4
$ws.Range("L107").Value = 1699.6666
$ws.Range("M107").Value = 310.5999999999999
$ws.Range("N107").Value = -5539.6666
$ws.Range("H113").Value = 482.33334
$ws.Range("I113").Value = 479
$ws.Range("J113").Value = 489
$ws.Range("K113").Value = 479
$ws.Range("L113").Value = 489
$ws.Range("M113").Value = 1691
$ws.Range("N113").Value = -4829
$ws.Range("H132").Value = 2158
$ws.Range("I132").Value = 2158
$ws.Range("K132").Value = 6474
$ws.Range("M132").Value = -3944
$ws.Range("H134").Value = 2525.3635
$ws.Range("I134").Value = 2707.389
$ws.Range("K134").Value = 8122.167
$ws.Range("M134").Value = -5587.167
$ws.Range("H135").Value = 110000
$ws.Range("J135").Value = 110000
$ws.Range("L135").Value = 110000
$ws.Range("N135").Value = -120140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4249.25
$ws.Range("I3").Value = 4249.25
$ws.Range("K3").Value = 12747.75
$ws.Range("M3").Value = -12635.75
$ws.Range("H58").Value = 11000
$ws.Range("I58").Value = 11000
$ws.Range("K58").Value = 33000
$ws.Range("M58").Value = -32872
$ws.Range("H121").Value = 1520927.5
$ws.Range("I121").Value = 166809.27
$ws.Range("J121").Value = 3777791
$ws.Range("K121").Value = 500427.8099999999
$ws.Range("L121").Value = 11333373
$ws.Range("M121").Value = -499117.8099999999
$ws.Range("N121").Value = -11335993
$ws.Range("H122").Value = 397.2
$ws.Range("I122").Value = 200
$ws.Range("J122").Value = 528.6667
$ws.Range("K122").Value = 1800
$ws.Range("L122").Value = 4758.0003
$ws.Range("M122").Value = 650
$ws.Range("N122").Value = -9658.0003
$ws.Range("H131").Value = 16777.527
$ws.Range("I131").Value = 371303.66
$ws.Range("J131").Value = 1797.5493
$ws.Range("K131").Value = 1113910.98
$ws.Range("L131").Value = 5392.6479
$ws.Range("M131").Value = -1108870.98
$ws.Range("N131").Value = -15472.6479
$ws.Range("H134").Value = 789.3333
$ws.Range("I134").Value = 789.3333
$ws.Range("K134").Value = 2367.9999
$ws.Range("M134").Value = 2702.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 39999
$ws.Range("J32").Value = 39999
$ws.Range("L32").Value = 39999
$ws.Range("N32").Value = -40591
$ws.Range("H80").Value = 14038.4
$ws.Range("I80").Value = 5239.7144
$ws.Range("K80").Value = 5239.7144
$ws.Range("M80").Value = -4241.7144
$ws.Range("H83").Value = 14038.4
$ws.Range("I83").Value = 5239.7144
$ws.Range("K83").Value = 26198.572
$ws.Range("M83").Value = -21206.572
$ws.Range("H126").Value = 2969.4
$ws.Range("I126").Value = 2565.5
$ws.Range("J126").Value = 3238.6667
$ws.Range("K126").Value = 7696.5
$ws.Range("L126").Value = 9716.000100000001
$ws.Range("M126").Value = -5226.5
$ws.Range("N126").Value = -14656.0001
$ws.Range("H132").Value = 2015.2
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H7").Value = 2429.1785
$ws.Range("J7").Value = 2938.2632
$ws.Range("L7").Value = 2938.2632
$ws.Range("N7").Value = -3162.2632
$ws.Range("H61").Value = 1914.5
$ws.Range("I61").Value = 1697.4
$ws.Range("K61").Value = 1697.4
$ws.Range("M61").Value = -1495.4
$ws.Range("H74").Value = 30000
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 30000
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H113").Value = 1914.5
$ws.Range("I113").Value = 1697.4
$ws.Range("K113").Value = 1697.4
$ws.Range("M113").Value = 472.5999999999999
$ws.Range("H126").Value = 2429.1785
$ws.Range("J126").Value = 2938.2632
$ws.Range("L126").Value = 8814.7896
$ws.Range("N126").Value = -13754.7896
$ws.Range("H132").Value = 4331.9165
$ws.Range("I132").Value = 3354.7144
$ws.Range("K132").Value = 10064.1432
$ws.Range("M132").Value = -7534.143199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 38349.668
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 38349.668
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 38349.668
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -39175.668
$ws.Range("H93").Value = 19999
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H126").Value = 1499
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
